$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: apply the "mtitleStyle" look (bold Century 12, centered, thin
# box border, black) to a label cell in the summary block (A10/A11/A12).
# ---------------------------------------------------------------------
function Set-MTitleStyle($cell) {
    $cell.Font.Name = "Century"
    $cell.Font.Size = 12
    $cell.Font.Bold = $true
    $cell.Font.Color = 0
    $cell.HorizontalAlignment = -4108
    $cell.Borders.LineStyle = 1
}

# Apply the "correctStyle" look (green Century 12, centered, thin box
# border) to a student-answer cell that was answered correctly.
function Set-CorrectStyle($cell) {
    $cell.Font.Name = "Century"
    $cell.Font.Size = 12
    $cell.Font.Bold = $false
    $cell.Font.Color = 32768
    $cell.HorizontalAlignment = -4108
    $cell.Borders.LineStyle = 1
}

# Apply the "incorrectStyle" look (red Century 12, centered, thin box
# border) to a student-answer cell that was answered incorrectly.
function Set-IncorrectStyle($cell) {
    $cell.Font.Name = "Century"
    $cell.Font.Size = 12
    $cell.Font.Bold = $false
    $cell.Font.Color = 255
    $cell.HorizontalAlignment = -4108
    $cell.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# Summary block (rows 10-12): the marking/float-parsing fix now produces
# real counts and a real score instead of all-zero / "Absent".
# ---------------------------------------------------------------------
Set-MTitleStyle $ws.Range("A10")
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 28

Set-MTitleStyle $ws.Range("A11")
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

Set-MTitleStyle $ws.Range("A12")
$ws.Range("B12").Value = 32
$ws.Range("C12").Value = -5
$ws.Range("E12").Value = "27/112"

# ---------------------------------------------------------------------
# Second answer block (columns D/E, "Student Ans" / "Correct Ans") only
# keeps its first three data rows (16-18) - and those now show the
# recorded (correct) student answers instead of being blank.
# ---------------------------------------------------------------------
Set-CorrectStyle $ws.Range("D16")
$ws.Range("D16").Value = "Option A"

Set-CorrectStyle $ws.Range("D17")
$ws.Range("D17").Value = "Option C"

Set-CorrectStyle $ws.Range("D18")
$ws.Range("D18").Value = "Option D"

# Drop the rest of the second block (rows 19-40) - it was only ever
# used for the unanswered placeholder rows.
$ws.Range("D19:E40").Clear()

# The whole third answer block (columns G/H) is removed entirely.
$ws.Range("G15:H21").Clear()

# ---------------------------------------------------------------------
# First answer block (columns A/B, "Student Ans" / "Correct Ans"):
# fill in the recorded student answers, colouring each green (correct)
# or red (incorrect) to match the grading in the Correct-Ans column.
# ---------------------------------------------------------------------
Set-CorrectStyle $ws.Range("A19")
$ws.Range("A19").Value = "Option C"

Set-IncorrectStyle $ws.Range("A20")
$ws.Range("A20").Value = "Option A"

Set-CorrectStyle $ws.Range("A21")
$ws.Range("A21").Value = "Option C"

Set-CorrectStyle $ws.Range("A24")
$ws.Range("A24").Value = "Option A"

Set-CorrectStyle $ws.Range("A26")
$ws.Range("A26").Value = "Option C"

Set-IncorrectStyle $ws.Range("A29")
$ws.Range("A29").Value = "Option C"

Set-IncorrectStyle $ws.Range("A30")
$ws.Range("A30").Value = "Option C"

Set-IncorrectStyle $ws.Range("A37")
$ws.Range("A37").Value = "Option B"

Set-CorrectStyle $ws.Range("A38")
$ws.Range("A38").Value = "Option A"

Set-IncorrectStyle $ws.Range("A39")
$ws.Range("A39").Value = "Option C"
